$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) in this sheet are stored as text
# (e.g. "314.87", "3.22%") rather than numbers/percentages. Pre-format the
# target range as Text so that assigning numeric-looking strings to
# .Value does not get auto-coerced into a real number by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "314.89"
$ws.Range("E2").Value = "3.62%"
$ws.Range("D3").Value = "35.10"
$ws.Range("E3").Value = "-1.31%"
$ws.Range("D4").Value = "5.096"
$ws.Range("E4").Value = "0.86%"
$ws.Range("D5").Value = "0.08160"
$ws.Range("E5").Value = "3.78%"
$ws.Range("D6").Value = "2.080"
$ws.Range("E6").Value = "-2.48%"
$ws.Range("D7").Value = "7.940"
$ws.Range("E7").Value = "-0.01%"
$ws.Range("D8").Value = "0.9308"
$ws.Range("E8").Value = "0.75%"
$ws.Range("D9").Value = "0.1038"
$ws.Range("D10").Value = "0.1927"
$ws.Range("E10").Value = "4.49%"
$ws.Range("D11").Value = "0.09115"
$ws.Range("E11").Value = "4.82%"
$ws.Range("D12").Value = "0.03600"
$ws.Range("E12").Value = "0.35%"
$ws.Range("D13").Value = "0.09902"
$ws.Range("E13").Value = "-0.40%"
$ws.Range("D14").Value = "0.001429"
$ws.Range("E14").Value = "-1.27%"
$ws.Range("D15").Value = "0.005714"
$ws.Range("E15").Value = "-0.23%"
$ws.Range("E16").Value = "-0.06%"
$ws.Range("D17").Value = "4.143"
$ws.Range("E17").Value = "0.17%"
$ws.Range("D18").Value = "2.853"
$ws.Range("E18").Value = "3.66%"
$ws.Range("E19").Value = "2.56%"
$ws.Range("D20").Value = "0.1292"
$ws.Range("E20").Value = "-4.16%"
$ws.Range("D21").Value = "5.101"
$ws.Range("E21").Value = "-1.39%"
$ws.Range("E22").Value = "0.10%"
$ws.Range("D23").Value = "0.04553"
$ws.Range("E23").Value = "-0.41%"
$ws.Range("E24").Value = "0.79%"
$ws.Range("D25").Value = "0.004791"
$ws.Range("E25").Value = "-0.71%"
$ws.Range("E26").Value = "-3.96%"
$ws.Range("D27").Value = "0.0004506"
$ws.Range("E27").Value = "-5.37%"
$ws.Range("D39").Value = "0.01979"
$ws.Range("E39").Value = "7.09%"
$ws.Range("D40").Value = "0.04949"
$ws.Range("E40").Value = "4.79%"
$ws.Range("D41").Value = "0.007580"
$ws.Range("E41").Value = "-2.75%"
$ws.Range("D42").Value = "0.1383"
$ws.Range("E42").Value = "-0.43%"
$ws.Range("D43").Value = "0.007875"
$ws.Range("E43").Value = "1.38%"
$ws.Range("D44").Value = "0.002253"
$ws.Range("E44").Value = "3.89%"
$ws.Range("D45").Value = "0.01176"
$ws.Range("E45").Value = "3.53%"
$ws.Range("D46").Value = "0.00006600"
$ws.Range("E46").Value = "4.74%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D48").Value = "188.28"
$ws.Range("E48").Value = "272.00%"
$ws.Range("D49").Value = "0.001702"
$ws.Range("E49").Value = "-10.64%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "-0.12%"
